$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell for the second column
$ws.Range("B1").Value = "AmountInCart"

# New checkout button/value cell
$ws.Range("B2").Value = "CHECKOUT ($539.98)"

# Column B sizing (auto/best-fit width similar to author's manual resize)
$ws.Columns.Item(2).ColumnWidth = 12.81640625

# Give column A (quantity cells) a top+bottom border instead of the old
# left/right+top+bottom "box" border, and move that original box border
# onto column B so the checkout button/amount cells keep it.
$rngA = $ws.Range("A2:A4")
$rngA.Borders.Item(8).LineStyle = 1
$rngA.Borders.Item(8).Weight = 2
$rngA.Borders.Item(8).Color = 0
$rngA.Borders.Item(9).LineStyle = 1
$rngA.Borders.Item(9).Weight = 2
$rngA.Borders.Item(9).Color = 0
$rngA.Borders.Item(7).LineStyle = -4142
$rngA.Borders.Item(10).LineStyle = -4142

$rngB = $ws.Range("B2:B4")
$rngB.Borders.Item(8).LineStyle = 1
$rngB.Borders.Item(8).Weight = 2
$rngB.Borders.Item(8).Color = 0
$rngB.Borders.Item(9).LineStyle = 1
$rngB.Borders.Item(9).Weight = 2
$rngB.Borders.Item(9).Color = 0
$rngB.Borders.Item(10).LineStyle = 1
$rngB.Borders.Item(10).Weight = 2
$rngB.Borders.Item(10).Color = 0
$rngB.Locked = $false
